$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 8 ----
$ws.Range("A8").Value = 10
$ws.Range("B8").Formula = "=A8*2"
$ws.Range("C8").Value = 100
$ws.Range("D8").Formula = "=A8/C8"
$ws.Range("E8").Formula = "=B8/C8"
$ws.Range("F8").Formula = "=(A8+B8)/(2*C8)"
$ws.Range("G8").Formula = "=(D8-E8) / ( SQRT( 2*F8*(1-F8)/C8 ) )"

# ---- Row 9 ----
$ws.Range("A9").Value = 100
$ws.Range("B9").Value = 200
$ws.Range("C9").Value = 1000
$ws.Range("D9").Formula = "=A9/C9"
$ws.Range("E9").Formula = "=B9/C9"
$ws.Range("F9").Formula = "=(A9+B9)/(2*C9)"
$ws.Range("G9").Formula = "=(D9-E9) / ( SQRT( 2*F9*(1-F9)/C9 ) )"

# ---- Row 10 ----
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 1000
$ws.Range("D10").Formula = "=A10/C10"
$ws.Range("E10").Formula = "=B10/C10"
$ws.Range("F10").Formula = "=(A10+B10)/(2*C10)"
$ws.Range("G10").Formula = "=(D10-E10) / ( SQRT( 2*F10*(1-F10)/C10 ) )"

# ---- Row 11 ----
$ws.Range("A11").Value = 60
$ws.Range("B11").Value = 30
$ws.Range("C11").Value = 1000
$ws.Range("D11").Formula = "=A11/C11"
$ws.Range("E11").Formula = "=B11/C11"
$ws.Range("F11").Formula = "=(A11+B11)/(2*C11)"
$ws.Range("G11").Formula = "=(D11-E11) / ( SQRT( 2*F11*(1-F11)/C11 ) )"

# ---- Row 12 ----
$ws.Range("A12").Value = 100
$ws.Range("B12").Value = 50
$ws.Range("C12").Value = 1000

# ---- Row 13 ----
$ws.Range("A13").Value = 100
$ws.Range("B13").Value = 95
$ws.Range("C13").Value = 1000

# D/E/F/G for rows 12:13 share formulas across the two-row range
$ws.Range("D12:D13").Formula = "=A12/C12"
$ws.Range("E12:E13").Formula = "=B12/C12"
$ws.Range("F12:F13").Formula = "=(A12+B12)/(2*C12)"
$ws.Range("G12:G13").Formula = "=(D12-E12) / ( SQRT( 2*F12*(1-F12)/C12 ) )"

Write-Output "data entered"

# ---- New note box H8:J9 ----
$noteStyle = $wb.Styles.Add("CenterWrapNote")
Write-Output "style created"
$noteStyle.HorizontalAlignment = -4108
$noteStyle.WrapText = $true
Write-Output "style configured"

$noteRng = $ws.Range("H8:J9")
$noteRng.Value = "More cells, but same proportions."
$noteRng.Style = $noteStyle
$noteRng.Merge()

Write-Output "note box done"

# ---- Selection ----
$ws.Range("N16").Select()

Write-Output "done"
